$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1 (base table rows 6-9) -> T-test summary in columns M/N
$ws.Range("M6").Value = "T-test"

$ws.Range("M7").Value = "2x2 depth 2"
$ws.Range("N7").Value = "86.3371 (3.6781)"

$ws.Range("M8").Value = "1x3 depth 2"
$ws.Range("N8").Value = "85.3271 (2.9752)"

$ws.Range("M9").Value = "P-value"
$ws.Range("N9").Value = 0.89004749999999999

# Block 2 (second table rows 12-16) -> T-test summary in columns M/N, offset by a row (row 11 onward)
$ws.Range("M11").Value = "T-test"

$ws.Range("M12").Value = "2x2 depth 2"
$ws.Range("N12").Value = "86.3371 (3.6781)"

$ws.Range("M13").Value = "1x2 depth 2"
$ws.Range("N13").Value = "85.0079 (2.8188)"

$ws.Range("M14").Value = "P-value"
$ws.Range("N14").Value = 0.84969519999999998

# Update selection to match final cursor position
$ws.Range("N19").Select()
